$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update realeffort (F column) scores for all ranked rows (2-13)
$ws.Range("F2").Value = 7.326165517881154
$ws.Range("F3").Value = 6.035906762210042
$ws.Range("F4").Value = 6.027117691378983
$ws.Range("F5").Value = 5.244593786151905
$ws.Range("F6").Value = 5.10415078822819
$ws.Range("F7").Value = 4.080915110249717
$ws.Range("F8").Value = 1.127328457611293
$ws.Range("F9").Value = 1.090887131911884

# Row 10: prolificid (B) and score (F) and race (G) updated
$ws.Range("B10").Value = 21
$ws.Range("F10").Value = 0.434968325099591
$ws.Range("G10").Value = "Black or African American"

# Row 11: prolificid (B) and score (F) updated
$ws.Range("B11").Value = 30
$ws.Range("F11").Value = 0.2911261319324809

# Row 12: score (F) updated
$ws.Range("F12").Value = 0.2519625011376062

# Row 13: prolificid (B), score (F) and race (G) updated
$ws.Range("B13").Value = 32
$ws.Range("F13").Value = 0.100228771449971
$ws.Range("G13").Value = "White"
